# Updated line loading-percentage results for the 380 kV case (Case_2_75).
# Values below come from res_line/loading_percent.xlsx, rows 2-25 (index 0-23),
# columns B:O, excluding D/I/L/N which remain 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "E", "F", "G", "H", "J", "K", "M", "O")

$newValues = @(
    @(8.07237213759527, 6.512864414584072, 22.09403190657481, 39.12605133402863, 26.56048421267059, 13.90256766727107, 8.007559211981116, 8.393281421867323, 17.54802599774226, 20.83846181071202),
    @(7.749888442166207, 6.450974607636262, 21.98251909181441, 39.1018885575188, 26.71961725655055, 13.95833118891351, 8.020243561036033, 8.152089365108889, 17.40116286893877, 20.94148841720897),
    @(7.545457142182801, 6.412815969858654, 21.91835575016834, 39.09715229159593, 26.82594073110578, 13.99466129540606, 8.02845537339339, 7.999041116220683, 17.31306760782924, 21.00904615021305),
    @(7.460652288691719, 6.397237151163668, 21.89331157538708, 39.09776444106825, 26.87142534170224, 14.00999252208824, 8.031908559026929, 7.935486393337694, 17.27772253212763, 21.03765694670859),
    @(7.446483796658128, 6.394648874626623, 21.88922017988474, 39.09801966868724, 26.87910802500005, 14.01257007243376, 8.032488417341966, 7.924863374582711, 17.2718879015235, 21.04247299056426),
    @(7.544319330832673, 6.412605970963395, 21.91801350461384, 39.09715025200465, 26.82654543158847, 13.99486592546768, 8.028501511419659, 7.998188715711542, 17.3125886457797, 21.00942763132252),
    @(7.962579876998153, 6.491564132924892, 22.05470112598596, 39.11562588711112, 26.61355955350428, 13.92136128297713, 8.011845041423625, 8.311181160843494, 17.49697527558681, 20.87309279951047),
    @(8.72730455506335, 6.644671526207248, 22.35595706765501, 39.23180111312269, 26.26467276258893, 13.79378355373219, 7.982529103854138, 8.883271754203442, 17.87362624143405, 20.63987555342926),
    @(9.250178083927327, 6.755474865346585, 22.59616506335212, 39.36551443638805, 26.05084901273964, 13.71011255047846, 7.963012390403673, 9.275349711846923, 18.15757844075192, 20.4893744171612),
    @(9.478724602086116, 6.805384763727071, 22.70922197085715, 39.43672348465761, 25.96292734409658, 13.67422522809042, 7.954568704086784, 9.447081107701552, 18.28789861916003, 20.42544191791515),
    @(9.563871658248404, 6.824202289875052, 22.75254945127564, 39.46516788355863, 25.93098811200776, 13.66094789199036, 7.951433487329205, 9.511122054316164, 18.3373765677607, 20.40188484843816),
    @(9.545596775074515, 6.820153431823726, 22.74319565724533, 39.45897633268266, 25.93780637417803, 13.66379351458074, 7.95210594913239, 9.497374269095751, 18.32671547336532, 20.40692923261631),
    @(9.485758041504187, 6.806934604365003, 22.71277639107897, 39.43903406598618, 25.96027247166823, 13.67312663436376, 7.954309522369983, 9.452369834193524, 18.29196683380786, 20.42349077389066),
    @(9.448921335366661, 6.798826652093746, 22.69420994080878, 39.42701102617331, 25.97421035281895, 13.6788841151265, 7.955667370156549, 9.424673358781044, 18.27069793883618, 20.43372024385686),
    @(9.23504937015373, 6.752202221819388, 22.58885023584554, 39.36106853342696, 26.05678390913766, 13.71250160796942, 7.963572926930318, 9.263990050336378, 18.14908169260043, 20.49364383711483),
    @(9.101415788118707, 6.723464890264622, 22.52516458376183, 39.32326483932441, 26.10984186443737, 13.73368167153279, 7.968533842427285, 9.163691075858736, 18.07474221841763, 20.53156640903898),
    @(9.023679835299122, 6.70688998621149, 22.48889235116981, 39.30249933157628, 26.14123862807862, 13.74606862839326, 7.971428148816187, 9.105379958092035, 18.03209445728169, 20.55380489445724),
    @(8.997211721631784, 6.701270457343584, 22.47667356374752, 39.29563688049475, 26.15201971811862, 13.75029781535408, 7.97241514724884, 9.085531217209635, 18.01767472707405, 20.5614076669026),
    @(9.115732169272897, 6.726528864944889, 22.53190717840033, 39.32718796358183, 26.10410268050919, 13.73140582994639, 7.968001511472171, 9.174432708329894, 18.082644628629, 20.52748534462261),
    @(9.503372516685591, 6.810819614331757, 22.72169751795731, 39.44485157165839, 25.95363678136327, 13.67037679472763, 7.953660592859932, 9.465615871668318, 18.30217015543627, 20.41860852167848),
    @(9.748544294073911, 6.865423991106828, 22.84872587615802, 39.53036492535219, 25.86320027442919, 13.63231163179064, 7.944650540948572, 9.650135970664682, 18.44637498967856, 20.35125672238821),
    @(9.618456898861126, 6.836328538907561, 22.7806647309626, 39.4839418020426, 25.910741454657, 13.65246122985157, 7.949426287495092, 9.552194665126921, 18.36935525174464, 20.38685503232029),
    @(9.109262557666728, 6.725143807576907, 22.52885778541971, 39.32541130174483, 26.10669458459255, 13.7324340829566, 7.968242046793987, 9.169578431527091, 18.07907166273029, 20.5293290348925),
    @(8.526943137262128, 6.603502514660818, 22.27104011512022, 39.19184503696123, 26.35163339955795, 13.82652734676365, 7.982529103854138, 8.883271754203442, 17.77031912355981, 20.69930931127683)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $addr = $cols[$j] + $rowNum
        $ws.Range($addr).Value = $rowData[$j]
    }
}

Write-Output "Updated loading_percent values for rows 2-25 (380 kV case)"
